$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold text-formatted figures (prices / percentages).
# Force text format so Excel does not reinterpret them as numbers and
# strip formatting (e.g. "42.711.01", "305.69", "  -0.86%  ").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.711.01"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.278.42"
$ws.Range("D5").Value = "305.69"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "96.65"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("D10").Value = "35.47"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "18.33"
$ws.Range("E12").Value = "  +3.23%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "2.634.71"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "2.286.36"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "42.615.60"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").Value = "67.09"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").Value = "235.93"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "25.10"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "165.70"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "33.01"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").Value = "17.59"
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "0.0690"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "2.70"
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("D43").Value = "1.997.79"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "18.16"
$ws.Range("E45").Value = "  +4.43%  "
$ws.Range("D46").Value = "9.95"
$ws.Range("E46").Value = "  -3.47%  "
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "53.64"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("D51").Value = "2.503.20"
$ws.Range("E51").Value = "  -1.10%  "
